$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09044833333333334
$ws.Range("H2").Value = 0.271345
$ws.Range("I2").Value = 0.2888886286400532
$ws.Range("J2").Value = 0.2888886286400532
$ws.Range("M2").Value = 201.5557555
$ws.Range("N2").Value = 403.111511
$ws.Range("O2").Value = 0.2814680640969941
$ws.Range("P2").Value = 0.2271476997658372
$ws.Range("Q2").Value = 18.23038215871583
$ws.Range("R2").Value = 109.382292952295
$ws.Range("S2").Value = 0.08131292304295122
$ws.Range("T2").Value = 0.06562038748409525
$ws.Range("G3").Value = 0.09044833333333334
$ws.Range("H3").Value = 0.271345
$ws.Range("I3").Value = 0.2888886286400532
$ws.Range("J3").Value = 0.2888886286400532
$ws.Range("N3").Value = 356.292984
$ws.Range("O3").Value = 0.1658517039268592
$ws.Range("P3").Value = 0.2007661144618275
$ws.Range("Q3").Value = 10.74203552705333
$ws.Range("R3").Value = 96.67831974348
$ws.Range("S3").Value = 0.04791267130504649
$ws.Range("T3").Value = 0.0579990474842693
$ws.Range("G4").Value = 0.09044833333333334
$ws.Range("H4").Value = 0.271345
$ws.Range("I4").Value = 0.2888886286400532
$ws.Range("J4").Value = 0.2888886286400532
$ws.Range("M4").Value = 77.61293766666667
$ws.Range("N4").Value = 232.838813
$ws.Range("O4").Value = 0.1083847159795808
$ws.Range("P4").Value = 0.1312014153551619
$ws.Range("Q4").Value = 7.019960857053889
$ws.Range("R4").Value = 63.179647713485
$ws.Range("S4").Value = 0.03131111196488277
$ws.Range("T4").Value = 0.03790259695758674
$ws.Range("G5").Value = 0.09044833333333334
$ws.Range("H5").Value = 0.271345
$ws.Range("I5").Value = 0.2888886286400532
$ws.Range("J5").Value = 0.2888886286400532
$ws.Range("M5").Value = 172.0397415
$ws.Range("N5").Value = 344.079483
$ws.Range("O5").Value = 0.2402496166265631
$ws.Range("P5").Value = 0.1938839774289365
$ws.Range("Q5").Value = 15.5607078857725
$ws.Range("R5").Value = 93.36424731463499
$ws.Range("S5").Value = 0.06940538227854634
$ws.Range("T5").Value = 0.05601087635472449
$ws.Range("G6").Value = 0.09044833333333334
$ws.Range("H6").Value = 0.271345
$ws.Range("I6").Value = 0.2888886286400532
$ws.Range("J6").Value = 0.2888886286400532
$ws.Range("M6").Value = 59.85226566666666
$ws.Range("N6").Value = 179.556797
$ws.Range("O6").Value = 0.08358233833226186
$ws.Range("P6").Value = 0.1011777443782085
$ws.Range("Q6").Value = 5.413537675773889
$ws.Range("R6").Value = 48.72183908196499
$ws.Range("S6").Value = 0.02414598709933608
$ws.Range("T6").Value = 0.0292290998223145
$ws.Range("G7").Value = 0.09044833333333334
$ws.Range("H7").Value = 0.271345
$ws.Range("I7").Value = 0.2888886286400532
$ws.Range("J7").Value = 0.2888886286400532
$ws.Range("M7").Value = 86.26244733333333
$ws.Range("N7").Value = 258.787342
$ws.Range("O7").Value = 0.1204635610377409
$ws.Range("P7").Value = 0.1458230486100285
$ws.Range("Q7").Value = 7.802294590554444
$ws.Range("R7").Value = 70.22065131498999
$ws.Range("S7").Value = 0.03480055294929033
$ws.Range("T7").Value = 0.04212662053706294
$ws.Range("G8").Value = 0.2226423333333333
$ws.Range("H8").Value = 0.6679269999999999
$ws.Range("I8").Value = 0.7111113713599468
$ws.Range("J8").Value = 0.7111113713599468
$ws.Range("M8").Value = 201.5557555
$ws.Range("N8").Value = 403.111511
$ws.Range("O8").Value = 0.2814680640969941
$ws.Range("P8").Value = 0.2271476997658372
$ws.Range("Q8").Value = 44.87484370128283
$ws.Range("R8").Value = 269.249062207697
$ws.Range("S8").Value = 0.2001551410540429
$ws.Range("T8").Value = 0.161527312281742
$ws.Range("G9").Value = 0.2226423333333333
$ws.Range("H9").Value = 0.6679269999999999
$ws.Range("I9").Value = 0.7111113713599468
$ws.Range("J9").Value = 0.7111113713599468
$ws.Range("N9").Value = 356.292984
$ws.Range("O9").Value = 0.1658517039268592
$ws.Range("P9").Value = 0.2007661144618275
$ws.Range("Q9").Value = 26.44196710268533
$ws.Range("R9").Value = 237.977703924168
$ws.Range("S9").Value = 0.1179390326218127
$ws.Range("T9").Value = 0.1427670669775582
$ws.Range("G10").Value = 0.2226423333333333
$ws.Range("H10").Value = 0.6679269999999999
$ws.Range("I10").Value = 0.7111113713599468
$ws.Range("J10").Value = 0.7111113713599468
$ws.Range("M10").Value = 77.61293766666667
$ws.Range("N10").Value = 232.838813
$ws.Range("O10").Value = 0.1083847159795808
$ws.Range("P10").Value = 0.1312014153551619
$ws.Range("Q10").Value = 17.27992553896122
$ws.Range("R10").Value = 155.519329850651
$ws.Range("S10").Value = 0.07707360401469808
$ws.Range("T10").Value = 0.09329881839757517
$ws.Range("G11").Value = 0.2226423333333333
$ws.Range("H11").Value = 0.6679269999999999
$ws.Range("I11").Value = 0.7111113713599468
$ws.Range("J11").Value = 0.7111113713599468
$ws.Range("M11").Value = 172.0397415
$ws.Range("N11").Value = 344.079483
$ws.Range("O11").Value = 0.2402496166265631
$ws.Range("P11").Value = 0.1938839774289365
$ws.Range("Q11").Value = 38.30332947362349
$ws.Range("R11").Value = 229.819976841741
$ws.Range("S11").Value = 0.1708442343480168
$ws.Range("T11").Value = 0.137873101074212
$ws.Range("G12").Value = 0.2226423333333333
$ws.Range("H12").Value = 0.6679269999999999
$ws.Range("I12").Value = 0.7111113713599468
$ws.Range("J12").Value = 0.7111113713599468
$ws.Range("M12").Value = 59.85226566666666
$ws.Range("N12").Value = 179.556797
$ws.Range("O12").Value = 0.08358233833226186
$ws.Range("P12").Value = 0.1011777443782085
$ws.Range("Q12").Value = 13.32564808331322
$ws.Range("R12").Value = 119.930832749819
$ws.Range("S12").Value = 0.05943635123292579
$ws.Range("T12").Value = 0.07194864455589399
$ws.Range("G13").Value = 0.2226423333333333
$ws.Range("H13").Value = 0.6679269999999999
$ws.Range("I13").Value = 0.7111113713599468
$ws.Range("J13").Value = 0.7111113713599468
$ws.Range("M13").Value = 86.26244733333333
$ws.Range("N13").Value = 258.787342
$ws.Range("O13").Value = 0.1204635610377409
$ws.Range("P13").Value = 0.1458230486100285
$ws.Range("Q13").Value = 19.20567255333711
$ws.Range("R13").Value = 172.851052980034
$ws.Range("S13").Value = 0.08566300808845062
$ws.Range("T13").Value = 0.1036964280729655
